$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Seed" column (D), shifting
# Seed/Train F1 Weighted/Test F1 Weighted/Y Val/Y Pred one column to the
# right (D->E, E->F, F->G, G->H, H->I).
$ws.Range("D1").EntireColumn.Insert()

# New header for the inserted column. The Insert() above already carried
# over the bold/bordered header style from the shifted columns, so we
# just need to set the text.
$ws.Range("D1").Value = "Best Score"

# Fill in the new "Best Score" values for each row.
$ws.Range("D2").Value = 0.4365681753771671
$ws.Range("D3").Value = 0.4102855230574304
$ws.Range("D4").Value = 0.4027467843474035
$ws.Range("D5").Value = 0.5359614936745103
$ws.Range("D6").Value = 0.4832580875773095

# Update the values of the shifted Train F1 Weighted (F) and
# Test F1 Weighted (G) columns, and refresh Y Val / Y Pred (H / I)
# arrays with their new values.
$ws.Range("F2").Value = 0.6535244922341696
$ws.Range("G2").Value = 0.5595238095238095
$ws.Range("H2").Value = "[1 0 1 0 0 1 1 1 1 1 1 1 1 0 1 0 0 0 1 0 1 1 0 0]"
$ws.Range("I2").Value = "[0 0 0 0 0 0 0 1 0 1 0 1 1 0 1 0 0 0 0 0 0 0 0 1]"

$ws.Range("F3").Value = 0.643601559730592
$ws.Range("G3").Value = 0.5370370370370371
$ws.Range("H3").Value = "[0 1 1 0 1 0 0 0 1 1 1 0 1 0 1 0 1 0 1 1 0 1 1 1]"
$ws.Range("I3").Value = "[0 1 1 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0]"

$ws.Range("F4").Value = 0.6120975202172421
$ws.Range("G4").Value = 0.5
$ws.Range("H4").Value = "[0 0 1 0 0 1 0 1 1 1 1 1 1 1 1 0 0 0 1 0 1 1 1 0]"
$ws.Range("I4").Value = "[0 1 1 0 0 0 1 1 1 0 0 1 1 0 0 0 1 0 0 1 1 0 0 0]"

$ws.Range("F5").Value = 0.6363079373832061
$ws.Range("G5").Value = 0.4140955837870539
$ws.Range("H5").Value = "[0 0 1 1 0 1 1 1 1 0 1 1 0 1 1 0 0 1 0 1 1 0 0 1]"
$ws.Range("I5").Value = "[0 0 1 0 0 0 1 0 0 1 0 0 1 0 1 0 0 0 0 0 0 0 0 0]"

$ws.Range("F6").Value = 0.7227056764191048
$ws.Range("G6").Value = 0.5
$ws.Range("H6").Value = "[1 0 1 0 1 1 0 0 1 1 0 1 0 1 1 1 1 1 0 1 0 0 1 0]"
$ws.Range("I6").Value = "[0 0 0 1 1 0 0 0 0 1 0 1 0 0 0 1 0 0 1 1 1 0 1 1]"
